# Edit LOM3233.xlsx per commit diff:
#  - Insert a new row (the sheet grows from 23 to 24 data rows) to correct
#    a misalignment in the "Objetivos/Programa resumido/Programa/Metodo/..."
#    block, and fill in several previously-missing texts (objectives text,
#    short syllabus summary, syllabus text, method text, criteria text,
#    recovery-norm text and the bibliography list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix rows 10/11 (Objetivos / Objectives): before this edit, B10/C10
#        incorrectly held the "519033 - Carlos Yujiro Shigue" text; put the
#        real Portuguese objectives text there instead.
$ws.Range("B10").Value = 'Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais.'
$ws.Range("C10").Value = 'Fornecer ao estudante noções básicas de dispositivos digitais e suas aplicações com ênfase em microcontroladores e processadores digitais de sinais.'

# --- 2. Insert a new row at 14 - this shifts old rows 14-23 down to 15-24,
#        and leaves a new blank row 14 (with inherited column styles).
$ws.Rows.Item(14).Insert()

# Row 13 ("Docentes responsáveis:" value row) should no longer have its own
# label in column A, and should lose its custom row height.
$ws.Range("A13").ClearContents()
$ws.Rows.Item(13).EntireRow.AutoFit()
$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'

# Row 14 becomes "Programa resumido:" with its short-syllabus text, height 60.
$ws.Rows.Item(14).RowHeight = 60
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle.'
$ws.Range("C14").Value = 'Circuitos digitais. Microprocessadores e microcontroladores. Programação de sistemas de aquisição de dados e algoritmos de controle.'

# --- 3. Row 16 ("Programa:") previously held a stray "01/01/2023" value;
#        replace it with the real (long) Portuguese syllabus text.
$ws.Range("B16").Value = "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA."
$ws.Range("C16").Value = "Bases numéricas. Aritmética binária. Funções lógicas. Álgebra de Boole. Minimização. Circuitos combinatórios. Flip-flops. Contadores e projeto de contadores. Introdução aos circuitos sequenciais. Microprocessadores. Microcontroladores e sistemas embarcados. Interfaces de comunicação. Linguagem de programação de baixo e alto nível na computação em tempo real. Desenvolvimento de protocolos de comando digital. Projeto com dispositivos programáveis: microcontroladores e processadores de sinais digitais. Programação de dispositivos FPGA."

# --- 4. Row 19 ("Método:") previously held the "519033 - Carlos Yujiro
#        Shigue" text; fix it to the real method text.
$ws.Range("B19").Value = 'Aulas expositivas, exercícios em sala, lista de exercícios, utilização de um simulador de circuitos, projeto de circuitos e atividades práticas em laboratório.'
$ws.Range("C19").Value = 'Aulas expositivas, exercícios em sala, lista de exercícios, utilização de um simulador de circuitos, projeto de circuitos e atividades práticas em laboratório.'

# --- 5. Row 20 ("Critério:") previously held the "Método" text; fix it to
#        the real evaluation-criteria text.
$ws.Range("B20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'
$ws.Range("C20").Value = 'Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4'

# --- 6. Row 21 ("Norma de recuperação:") previously held the "Critério"
#        text; fix it to the real recovery-norm text.
$ws.Range("B21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

# --- 7. Row 22 ("Bibliografia:") previously held the "Norma de
#        recuperação" text; fix it to the real multi-line bibliography list.
$bibliografia = "GAJSKI, D. D. Principles of Digital Design, Prentice Hall, 1997.`nTAUB, H. Circuitos Digitais e Microprocessadores, McGraw Hill, 1984.`nTOCCI, R. J.; AMBROSIO, F. J. Microprocessors and Microcomputers: Hardware and Software, Prentice Hall, 2002.`nCATSOULIS, J. Designing Embedded Hardware, OReilly Media, 2005.`nCRISP, J. Introduction to Microprocessors, Newnes, 2004.`nWILMSHURST, T. Designing Embedded Systems with PIC Microcontrollers, Newnes, 2009.`nDUBEY, R. Introduction to Embedded System Design using Field Programmable Gate Arrays, Springer, 2008.`nBATEMAN, A.; PATERSON-STEPHENS, I. The DSP Handbook: Algorithms, Applications and Design Techniques, Prentice Hall, 2002."
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# --- 8. Tidy the column definitions: column A's width/style no longer
#        needs to be declared jointly with column B (split the merged
#        min=1,max=2 column range into its own min=1,max=1 entry).
$ws.Columns.Item(2).ColumnWidth = 60.7109375
